# Update the dSF column (F) values with repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -11
    3  = -6
    4  = -2
    5  = 4
    6  = -3
    7  = 3
    8  = -2
    9  = -5
    10 = 8
    11 = 4
    12 = -3
    13 = 1
    14 = -1
    15 = -2
    16 = -1
    17 = 5
    18 = 0
    19 = -2
    20 = -4
    21 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
